$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel
# (e.g. "213.00" -> 213) need to be pre-formatted as Text so the literal
# string is preserved, matching the source data export.
$textCells = @("D5", "D6", "D8", "D9", "D15", "D18", "D19", "D20", "D23", "D24", "D26", "D27", "D28", "D31", "D38", "D40", "D44", "D45", "D47", "D48", "D50")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.552.65"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "1.570.77"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "213.00"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Value = "0.491"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "45.67"
$ws.Range("E8").Value = "  +4.17%  "
$ws.Range("D9").Value = "24.03"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "1.794.52"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").Value = "1.570.66"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Value = "0.521"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "28.548.25"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "62.28"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").Value = "230.18"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").Value = "7.36"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("E21").Value = "  -2.71%  "
$ws.Range("D23").Value = "3.87"
$ws.Range("E23").Value = "  -5.64%  "
$ws.Range("D24").Value = "9.12"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("E25").Value = "  +9.60%  "
$ws.Range("D26").Value = "151.82"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "15.03"
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").Value = "6.43"
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("E29").Value = "  -3.28%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "0.0486"
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("D35").Value = "1.394.54"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("E37").Value = "  -3.28%  "
$ws.Range("D38").Value = "2.37"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("E39").Value = "  +2.75%  "
$ws.Range("D40").Value = "0.0166"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").Value = "0.790"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("D45").Value = "0.0470"
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").Value = "0.969"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").Value = "62.92"
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("D49").Value = "1.707.04"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").Value = "86.34"
$ws.Range("E51").Value = "  -0.87%  "
